$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay stored as text, matching the
# original inline-string cells in the workbook (avoids Excel auto-converting
# values like "39.94" or "0.600" into floating point numbers).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.132.66"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "3.823.82"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "420.75"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "127.34"
$ws.Range("E6").Value = "  -4.25%  "
$ws.Range("D7").Value = "3.818.63"
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  -7.96%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -8.05%  "
$ws.Range("E11").Value = "  -14.23%  "
$ws.Range("D12").Value = "0.0000343"
$ws.Range("E12").Value = "  -20.30%  "
$ws.Range("D13").Value = "39.94"
$ws.Range("E13").Value = "  -7.20%  "
$ws.Range("D14").Value = "4.435.58"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").Value = "15.89"
$ws.Range("E15").Value = "  +21.92%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "9.85"
$ws.Range("E16").Value = "  -5.87%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.137"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.809.11"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "19.37"
$ws.Range("E19").Value = "  -6.13%  "
$ws.Range("D20").Value = "66.214.59"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("E21").Value = "  -7.06%  "
$ws.Range("D22").Value = "400.04"
$ws.Range("E22").Value = "  -11.39%  "
$ws.Range("E23").Value = "  -11.22%  "
$ws.Range("D24").Value = "83.44"
$ws.Range("E24").Value = "  -8.15%  "
$ws.Range("E25").Value = "  -5.02%  "
$ws.Range("D26").Value = "36.88"
$ws.Range("E26").Value = "  -4.76%  "
$ws.Range("E27").Value = "  +12.53%  "
$ws.Range("D28").Value = "3.16"
$ws.Range("E28").Value = "  -5.99%  "
$ws.Range("D29").Value = "9.31"
$ws.Range("E29").Value = "  -8.19%  "
$ws.Range("D30").Value = "696.34"
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("D33").Value = "12.18"
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("D34").Value = "7.41"
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("E35").Value = "  -10.70%  "
$ws.Range("D36").Value = "37.49"
$ws.Range("E36").Value = "  -10.80%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "54.64"
$ws.Range("E38").Value = "  -5.96%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "0.0448"
$ws.Range("E40").Value = "  -9.32%  "
$ws.Range("D41").Value = "2.89"
$ws.Range("E41").Value = "  -4.54%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "0.134"
$ws.Range("E43").Value = "  -10.17%  "
$ws.Range("D44").Value = "4.41"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").Value = "3.29"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("D46").Value = "142.89"
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "2.04"
$ws.Range("E47").Value = "  -4.48%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.05"
$ws.Range("E48").Value = "  -3.92%  "
$ws.Range("D49").Value = "25.72"
$ws.Range("E49").Value = "  -7.64%  "
$ws.Range("D50").Value = "2.51"
$ws.Range("E50").Value = "  -5.19%  "
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").Value = "  -7.47%  "
